$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 3.537761648806719
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 1138.349488340621

$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 6.189590430959694

$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 6708.013860684405
$ws.Range("E4").Value = 10.19245300693656
$ws.Range("G4").Value = 6723.148924318466

$ws.Range("B5").Value = 0.0001021024915524027
$ws.Range("C5").Value = 0.00006240767534437808
$ws.Range("D5").Value = 0.7527432677738641
$ws.Range("E5").Value = 10.19245300693656
$ws.Range("G5").Value = 10.94536078487732

$ws.Range("B6").Value = 3.286832544864788
$ws.Range("C6").Value = 1.655778082260271
$ws.Range("D6").Value = 3.537761648806719
$ws.Range("E6").Value = 0.4942365360607697
$ws.Range("G6").Value = 8.974608811992548
